# Regenerate save_data: update column G ("K", formerly "Strike#") values
# for rows 2-34 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G (K)
$newValues = @{
    2  = 5
    3  = 1
    4  = 3
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 3
    13 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 4
    19 = 3
    20 = 9
    21 = 2
    22 = 4
    23 = 4
    24 = 5
    25 = 5
    26 = 4
    27 = 1
    28 = 3
    29 = 5
    30 = 3
    31 = 10
    32 = 6
    33 = 3
    34 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
